# Re-apply the "28/7/2025(Remote)" work-log entry that had been lost, adding it
# back as a new row (row 21) at the bottom of the Car Tracking Project log on
# Sheet1 -- matching the columns Date / Project / Task / Context / Problem / Status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date / Project / Task / Problem / Status for the new entry (no Context/D value).
$ws.Cells.Item(21, 1).Value = "28/7/2025(Remote)"
$ws.Cells.Item(21, 2).Value = "Car Tracking Project"
$ws.Cells.Item(21, 3).Value = "Try to make the prefect and the github actions to run the scripts daily without my need to run it `nmanually"
$ws.Cells.Item(21, 5).Value = "While changing my code on a branch, the merged it into a the main branch by a mistake; therefore all of todays progression is gone. And I will have to `nmodularize the code again."
$ws.Cells.Item(21, 6).Value = "FAILED REALLY HARD"

# Task (C) and Problem (E) entries wrap, like the other multi-line rows above them.
$ws.Range("C21").WrapText = $true
$ws.Range("E21").WrapText = $true

# Two wrapped lines -> same 28.8pt row height used elsewhere in the sheet for 2-line rows.
$ws.Rows.Item(21).RowHeight = 28.8

# Scroll the view down/over towards the newly added row and reselect F21, as last left by the author.
try {
    $excel.ActiveWindow.ScrollRow = 17
    $excel.ActiveWindow.ScrollColumn = 5
} catch {
}
$ws.Range("F21").Select() | Out-Null
